$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.153.22"

$ws.Range("D3").Value = "3.559.18"
$ws.Range("E3").Value = "  +1.99%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.56"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.33"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").Value = "3.557.28"
$ws.Range("E7").Value = "  +1.97%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  +3.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.82"
$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "4.165.42"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000206"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.97"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "3.575.94"
$ws.Range("E16").Value = "  +2.57%  "

$ws.Range("D17").Value = "66.235.76"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  +5.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.12"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.86"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").Value = "3.706.09"
$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.10"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.84"
$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").Value = "3.559.95"
$ws.Range("E32").Value = "  +2.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.42"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("E35").Value = "  -7.83%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.79"
$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.71"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "173.10"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.17"
$ws.Range("E42").Value = "  +0.27%  "

$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.04"
$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.97"
$ws.Range("E48").Value = "  -3.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.40"
$ws.Range("E49").Value = "  +0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.99"
$ws.Range("E51").Value = "  +4.10%  "
